$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44461
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112026
$ws.Cells.Item($row, 7).Value = "Haba"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 1100
$ws.Cells.Item($row, 11).Value = 9000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 9500
$ws.Cells.Item($row, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 380
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
